$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 01_리그테이블
# Insert a new row for "디비금융스팩12호" (DB) right after the 라메디텍 row,
# i.e. as the new row 3, pushing every following row down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("01_리그테이블")
$ws1.Rows.Item(3).Insert()

$row1 = $ws1.Rows.Item(3)
$row1.Cells.Item(1, 1).Value = "2024-06-18"
$row1.Cells.Item(1, 2).Value = "디비금융스팩12호"
$row1.Cells.Item(1, 3).Value = "코스닥"
$row1.Cells.Item(1, 4).Value = 100
$row1.Cells.Item(1, 5).Value = "DB"
$row1.Cells.Item(1, 6).Value = 100
$row1.Cells.Item(1, 7).Value = "-"
$row1.Cells.Item(1, 8).Value = "-"
$row1.Cells.Item(1, 9).Value = "-"
$row1.Cells.Item(1, 10).Value = "-"
$row1.Cells.Item(1, 11).Value = "대표"
$row1.Cells.Item(1, 12).Value = "-"
$row1.Cells.Item(1, 13).Value = 2000
$row1.Cells.Item(1, 14).Value = 100
$row1.Cells.Item(1, 15).Value = "2024-06-05"
$row1.Cells.Item(1, 16).Value = "2024-06-11"
$row1.Cells.Item(1, 17).Value = 3750000

# ---------------------------------------------------------------------------
# Sheet 2: 02_통합집계_Rawdata
# Insert the matching raw-data row, also right after the 라메디텍 row
# (new row 3).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("02_통합집계_Rawdata")
$ws2.Rows.Item(3).Insert()

$row2 = $ws2.Rows.Item(3)
$row2.Cells.Item(1, 1).Value = "2024-06-05"
$row2.Cells.Item(1, 2).Value = "디비금융스팩12호"
$row2.Cells.Item(1, 3).Value = "DB"
$row2.Cells.Item(1, 4).Value = "2024-06-11"
$row2.Cells.Item(1, 5).Value = "2024-06-18"
$row2.Cells.Item(1, 6).Value = 10000000
$row2.Cells.Item(1, 7).Value = 5000000
$row2.Cells.Item(1, 8).Value = "-"
$row2.Cells.Item(1, 9).Value = 2000
$row2.Cells.Item(1, 10).Value = 2000
$row2.Cells.Item(1, 11).Value = "-"
$row2.Cells.Item(1, 12).Value = 2000
$row2.Cells.Item(1, 13).Value = "-"
$row2.Cells.Item(1, 14).Value = "-"
$row2.Cells.Item(1, 15).Value = 0
$row2.Cells.Item(1, 16).Value = "-"
$row2.Cells.Item(1, 17).Value = "-"
$row2.Cells.Item(1, 18).Value = "1141.4 : 1"
$row2.Cells.Item(1, 19).Value = "-"
$row2.Cells.Item(1, 20).Value = "-"

# ---------------------------------------------------------------------------
# Sheet 3: 03_IPO현황_Summary
# This sheet is sorted alphabetically by 인수기관 (underwriter). "DB" sorts
# before "KB", so the new row becomes the new row 2.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("03_IPO현황_Summary")
$ws3.Rows.Item(2).Insert()

$row3 = $ws3.Rows.Item(2)
$row3.Cells.Item(1, 1).Value = "DB"
$row3.Cells.Item(1, 2).Value = "2024-06-05"
$row3.Cells.Item(1, 3).Value = "디비금융스팩12호"
$row3.Cells.Item(1, 4).Value = "DB"
$row3.Cells.Item(1, 5).Value = "DB"
$row3.Cells.Item(1, 6).Value = "2024-06-11"
$row3.Cells.Item(1, 7).Value = "2024-06-18"
$row3.Cells.Item(1, 8).Value = 10000
$row3.Cells.Item(1, 9).Value = 5000000
$row3.Cells.Item(1, 10).Value = 2000
$row3.Cells.Item(1, 11).Value = 0
$row3.Cells.Item(1, 12).Value = 100
